$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.63"
$ws.Range("E2").Value = "'-0.60%"
$ws.Range("D3").Value = "'37.46"
$ws.Range("E3").Value = "'7.46%"
$ws.Range("E4").Value = "'-2.48%"
$ws.Range("D5").Value = "'0.07863"
$ws.Range("E5").Value = "'1.21%"
$ws.Range("D6").Value = "'2.252"
$ws.Range("E6").Value = "'-4.98%"
$ws.Range("D7").Value = "'8.023"
$ws.Range("E7").Value = "'0.02%"
$ws.Range("D8").Value = "'4.018"
$ws.Range("E8").Value = "'1.99%"
$ws.Range("D9").Value = "'0.9098"
$ws.Range("E9").Value = "'-1.89%"
$ws.Range("D10").Value = "'0.1876"
$ws.Range("E10").Value = "'4.50%"
$ws.Range("D11").Value = "'0.09216"
$ws.Range("E11").Value = "'-9.25%"
$ws.Range("D12").Value = "'0.08523"
$ws.Range("E12").Value = "'-0.42%"
$ws.Range("D13").Value = "'0.03533"
$ws.Range("E13").Value = "'6.69%"
$ws.Range("D14").Value = "'0.09938"
$ws.Range("E14").Value = "'0.42%"
$ws.Range("D15").Value = "'0.001491"
$ws.Range("E15").Value = "'-0.75%"
$ws.Range("D16").Value = "'0.005623"
$ws.Range("E16").Value = "'-2.17%"
$ws.Range("E17").Value = "'-0.01%"
$ws.Range("D18").Value = "'2.157"
$ws.Range("E18").Value = "'-0.37%"
$ws.Range("E19").Value = "'2.86%"
$ws.Range("E20").Value = "'-0.20%"
$ws.Range("D21").Value = "'4.745"
$ws.Range("E21").Value = "'10.58%"
$ws.Range("E22").Value = "'-7.83%"
$ws.Range("D23").Value = "'0.04651"
$ws.Range("E23").Value = "'1.90%"
$ws.Range("E24").Value = "'0.85%"
$ws.Range("D25").Value = "'0.004451"
$ws.Range("E25").Value = "'-0.27%"
$ws.Range("E26").Value = "'3.81%"
$ws.Range("E27").Value = "'28.15%"
$ws.Range("D39").Value = "'0.01760"
$ws.Range("E39").Value = "'-1.50%"
$ws.Range("D40").Value = "'0.04733"
$ws.Range("E40").Value = "'-0.61%"
$ws.Range("D41").Value = "'0.007855"
$ws.Range("E41").Value = "'1.56%"
$ws.Range("E42").Value = "'-1.47%"
$ws.Range("D43").Value = "'0.007650"
$ws.Range("E43").Value = "'8.05%"
$ws.Range("E44").Value = "'5.25%"
$ws.Range("D45").Value = "'0.01020"
$ws.Range("E45").Value = "'6.97%"
$ws.Range("D46").Value = "'0.00005977"
$ws.Range("E46").Value = "'-2.32%"
$ws.Range("E47").Value = "'-0.17%"
$ws.Range("E48").Value = "'216.96%"
$ws.Range("E49").Value = "'34.26%"
$ws.Range("E50").Value = "'-0.17%"
$ws.Range("D51").Value = "'0.0001996"
$ws.Range("E51").Value = "'-0.17%"
